# UT - Data Collection Types: add a second "#" record (index 2) to the
# array/list/set comparison table, and make the "list" collection for
# that second record shorter (2 elements instead of 3) so the reader
# has to cope with differing collection sizes across records.
#
# Layout recap (columns C/D = array, E/F = list, G/H = set; each pair is
# boolean Property / text Property). Record 1 already occupies rows 4-6.
# Record 2 now occupies rows 7-9, but its "list" column only has 2
# elements (row 9, cols E/F stay empty) instead of 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UT - Data Collection Types")

# --- Record 1 (index 1, rows 4-6): trim its 3rd "array"/"set" element
#     (row 6) down to just the "list" pair so the 3 collections still
#     align under the new layout; array/set 3rd elements move to the
#     new record-2 rows below.
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()

# --- Record 2 (index 2): row 7 = "#" index + first element of each
#     collection (array/list/set all populated, 3 elements each).
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "'true"
$ws.Range("D7").Value = "'ABC"
$ws.Range("E7").Value = "'true"
$ws.Range("F7").Value = "'ABC"
$ws.Range("G7").Value = "'true"
$ws.Range("H7").Value = "'ABC"

# Row 8 = second element of each collection for record 2.
$ws.Range("C8").Value = "'false"
$ws.Range("D8").Value = "'BCD"
$ws.Range("E8").Value = "'false"
$ws.Range("F8").Value = "'BCD"
$ws.Range("G8").Value = "'false"
$ws.Range("H8").Value = "'BCD"

# Row 9 = third element of "array" and "set" only; "list" deliberately
# stops at 2 elements for this record, so E9/F9 stay blank.
$ws.Range("C9").Value = "'true"
$ws.Range("D9").Value = "'ABC"
$ws.Range("G9").Value = "'true"
$ws.Range("H9").Value = "'ABC"

# Match the formatting used by the equivalent cells of record 1.
$ws.Range("B7").HorizontalAlignment = -4108
$ws.Range("B7").VerticalAlignment = -4108

$ws.Range("C8").HorizontalAlignment = -4131
$ws.Range("C8").VerticalAlignment = -4108
$ws.Range("E8").HorizontalAlignment = -4131
$ws.Range("E8").VerticalAlignment = -4108
$ws.Range("G8").HorizontalAlignment = -4131
$ws.Range("G8").VerticalAlignment = -4108

# Leave the selection where the author left it.
$ws.Activate()
$ws.Range("H9").Select()
